$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure text columns (League, Date, Time, Home, Away) stay as plain text,
# so Excel does not auto-convert date/time-like strings into date serials.
$ws.Range("A2:E14").NumberFormat = "@"

# Row 2: Egyptian Premier | Al Ahly Cairo vs Wadi Degla
$ws.Cells.Item(2, 1).Value = "Egyptian Premier"
$ws.Cells.Item(2, 2).Value = "2026-01-27"
$ws.Cells.Item(2, 3).Value = "12:00:00"
$ws.Cells.Item(2, 4).Value = "Al Ahly Cairo"
$ws.Cells.Item(2, 5).Value = "Wadi Degla"
$ws.Cells.Item(2, 6).Value = 1.49
$ws.Cells.Item(2, 7).Value = 1.57
$ws.Cells.Item(2, 8).Value = 8
$ws.Cells.Item(2, 9).Value = 10.5
$ws.Cells.Item(2, 10).Value = 4.1
$ws.Cells.Item(2, 11).Value = 4.6
$ws.Cells.Item(2, 12).Value = 1.01
$ws.Cells.Item(2, 13).Value = 1.08
$ws.Cells.Item(2, 14).Value = 3.2
$ws.Cells.Item(2, 15).Value = 1.38
$ws.Cells.Item(2, 16).Value = 1.75
$ws.Cells.Item(2, 17).Value = 2.14
$ws.Cells.Item(2, 18).Value = 1.28
$ws.Cells.Item(2, 19).Value = 3.85
$ws.Cells.Item(2, 20).Value = 2.06
$ws.Cells.Item(2, 21).Value = 1.64
$ws.Cells.Item(2, 22).Value = 1.1
$ws.Cells.Item(2, 23).Value = 2.74
$ws.Cells.Item(2, 24).Value = 14
$ws.Cells.Item(2, 25).Value = 23
$ws.Cells.Item(2, 26).Value = 80
$ws.Cells.Item(2, 27).Value = 1000
$ws.Cells.Item(2, 28).Value = 6.8
$ws.Cells.Item(2, 29).Value = 10.5
$ws.Cells.Item(2, 30).Value = 36
$ws.Cells.Item(2, 31).Value = 210
$ws.Cells.Item(2, 32).Value = 8
$ws.Cells.Item(2, 33).Value = 11.5
$ws.Cells.Item(2, 34).Value = 34
$ws.Cells.Item(2, 35).Value = 190
$ws.Cells.Item(2, 36).Value = 14.5
$ws.Cells.Item(2, 37).Value = 21
$ws.Cells.Item(2, 38).Value = 55
$ws.Cells.Item(2, 39).Value = 270
$ws.Cells.Item(2, 40).Value = 11
$ws.Cells.Item(2, 41).Value = 360

# Row 3: Dutch Eerste Divisie | Emmen vs Cambuur Leeuwarden
$ws.Cells.Item(3, 1).Value = "Dutch Eerste Divisie"
$ws.Cells.Item(3, 2).Value = "2026-01-27"
$ws.Cells.Item(3, 3).Value = "16:00:00"
$ws.Cells.Item(3, 4).Value = "Emmen"
$ws.Cells.Item(3, 5).Value = "Cambuur Leeuwarden"
$ws.Cells.Item(3, 6).Value = 3.8
$ws.Cells.Item(3, 7).Value = 4.5
$ws.Cells.Item(3, 8).Value = 1.86
$ws.Cells.Item(3, 9).Value = 2
$ws.Cells.Item(3, 10).Value = 3.95
$ws.Cells.Item(3, 11).Value = 4.4
$ws.Cells.Item(3, 12).Value = 1.01
$ws.Cells.Item(3, 13).Value = 1.03
$ws.Cells.Item(3, 14).Value = 5.5
$ws.Cells.Item(3, 15).Value = 1.18
$ws.Cells.Item(3, 16).Value = 2.54
$ws.Cells.Item(3, 17).Value = 1.55
$ws.Cells.Item(3, 18).Value = 1.57
$ws.Cells.Item(3, 19).Value = 2.2
$ws.Cells.Item(3, 20).Value = 1.53
$ws.Cells.Item(3, 21).Value = 2.54
$ws.Cells.Item(3, 22).Value = 2
$ws.Cells.Item(3, 23).Value = 1.29
$ws.Cells.Item(3, 24).Value = 26
$ws.Cells.Item(3, 25).Value = 14
$ws.Cells.Item(3, 26).Value = 16
$ws.Cells.Item(3, 27).Value = 24
$ws.Cells.Item(3, 28).Value = 25
$ws.Cells.Item(3, 29).Value = 10.5
$ws.Cells.Item(3, 30).Value = 11
$ws.Cells.Item(3, 31).Value = 18.5
$ws.Cells.Item(3, 32).Value = 38
$ws.Cells.Item(3, 33).Value = 19.5
$ws.Cells.Item(3, 34).Value = 16
$ws.Cells.Item(3, 35).Value = 27
$ws.Cells.Item(3, 36).Value = 80
$ws.Cells.Item(3, 37).Value = 44
$ws.Cells.Item(3, 38).Value = 46
$ws.Cells.Item(3, 39).Value = 60
$ws.Cells.Item(3, 40).Value = 32
$ws.Cells.Item(3, 41).Value = 9.2

# Row 4: German Bundesliga | St Pauli vs RB Leipzig
$ws.Cells.Item(4, 1).Value = "German Bundesliga"
$ws.Cells.Item(4, 2).Value = "2026-01-27"
$ws.Cells.Item(4, 3).Value = "16:30:00"
$ws.Cells.Item(4, 4).Value = "St Pauli"
$ws.Cells.Item(4, 5).Value = "RB Leipzig"
$ws.Cells.Item(4, 6).Value = 4.4
$ws.Cells.Item(4, 7).Value = 4.6
$ws.Cells.Item(4, 8).Value = 1.89
$ws.Cells.Item(4, 9).Value = 1.91
$ws.Cells.Item(4, 10).Value = 3.9
$ws.Cells.Item(4, 11).Value = 3.95
$ws.Cells.Item(4, 12).Value = 1.35
$ws.Cells.Item(4, 13).Value = 1.05
$ws.Cells.Item(4, 14).Value = 4.3
$ws.Cells.Item(4, 15).Value = 1.27
$ws.Cells.Item(4, 16).Value = 2.18
$ws.Cells.Item(4, 17).Value = 1.81
$ws.Cells.Item(4, 18).Value = 1.45
$ws.Cells.Item(4, 19).Value = 3.05
$ws.Cells.Item(4, 20).Value = 1.76
$ws.Cells.Item(4, 21).Value = 2.24
$ws.Cells.Item(4, 22).Value = 2.08
$ws.Cells.Item(4, 23).Value = 1.27
$ws.Cells.Item(4, 24).Value = 18
$ws.Cells.Item(4, 25).Value = 10.5
$ws.Cells.Item(4, 26).Value = 12
$ws.Cells.Item(4, 27).Value = 21
$ws.Cells.Item(4, 28).Value = 19
$ws.Cells.Item(4, 29).Value = 8.6
$ws.Cells.Item(4, 30).Value = 10.5
$ws.Cells.Item(4, 31).Value = 18
$ws.Cells.Item(4, 32).Value = 34
$ws.Cells.Item(4, 33).Value = 18.5
$ws.Cells.Item(4, 34).Value = 17.5
$ws.Cells.Item(4, 35).Value = 32
$ws.Cells.Item(4, 36).Value = 100
$ws.Cells.Item(4, 37).Value = 55
$ws.Cells.Item(4, 38).Value = 60
$ws.Cells.Item(4, 39).Value = 90
$ws.Cells.Item(4, 40).Value = 48
$ws.Cells.Item(4, 41).Value = 11

# Row 5: German Bundesliga | Werder Bremen vs Hoffenheim
$ws.Cells.Item(5, 1).Value = "German Bundesliga"
$ws.Cells.Item(5, 2).Value = "2026-01-27"
$ws.Cells.Item(5, 3).Value = "16:30:00"
$ws.Cells.Item(5, 4).Value = "Werder Bremen"
$ws.Cells.Item(5, 5).Value = "Hoffenheim"
$ws.Cells.Item(5, 6).Value = 3.25
$ws.Cells.Item(5, 7).Value = 3.35
$ws.Cells.Item(5, 8).Value = 2.28
$ws.Cells.Item(5, 9).Value = 2.32
$ws.Cells.Item(5, 10).Value = 3.85
$ws.Cells.Item(5, 11).Value = 3.9
$ws.Cells.Item(5, 12).Value = 1.28
$ws.Cells.Item(5, 13).Value = 1.04
$ws.Cells.Item(5, 14).Value = 5.4
$ws.Cells.Item(5, 15).Value = 1.2
$ws.Cells.Item(5, 16).Value = 2.46
$ws.Cells.Item(5, 17).Value = 1.63
$ws.Cells.Item(5, 18).Value = 1.6
$ws.Cells.Item(5, 19).Value = 2.56
$ws.Cells.Item(5, 20).Value = 1.56
$ws.Cells.Item(5, 21).Value = 2.68
$ws.Cells.Item(5, 22).Value = 1.76
$ws.Cells.Item(5, 23).Value = 1.43
$ws.Cells.Item(5, 24).Value = 23
$ws.Cells.Item(5, 25).Value = 14.5
$ws.Cells.Item(5, 26).Value = 17
$ws.Cells.Item(5, 27).Value = 30
$ws.Cells.Item(5, 28).Value = 18.5
$ws.Cells.Item(5, 29).Value = 9
$ws.Cells.Item(5, 30).Value = 11.5
$ws.Cells.Item(5, 31).Value = 21
$ws.Cells.Item(5, 32).Value = 26
$ws.Cells.Item(5, 33).Value = 13.5
$ws.Cells.Item(5, 34).Value = 14.5
$ws.Cells.Item(5, 35).Value = 29
$ws.Cells.Item(5, 36).Value = 55
$ws.Cells.Item(5, 37).Value = 32
$ws.Cells.Item(5, 38).Value = 36
$ws.Cells.Item(5, 39).Value = 60
$ws.Cells.Item(5, 40).Value = 22
$ws.Cells.Item(5, 41).Value = 12.5

# Row 6: English National League | Rochdale vs Southend
$ws.Cells.Item(6, 1).Value = "English National League"
$ws.Cells.Item(6, 2).Value = "2026-01-27"
$ws.Cells.Item(6, 3).Value = "16:45:00"
$ws.Cells.Item(6, 4).Value = "Rochdale"
$ws.Cells.Item(6, 5).Value = "Southend"
$ws.Cells.Item(6, 6).Value = 2.18
$ws.Cells.Item(6, 7).Value = 2.4
$ws.Cells.Item(6, 8).Value = 3.15
$ws.Cells.Item(6, 9).Value = 3.65
$ws.Cells.Item(6, 10).Value = 3.55
$ws.Cells.Item(6, 11).Value = 4.1
$ws.Cells.Item(6, 12).Value = 1.01
$ws.Cells.Item(6, 13).Value = 1.05
$ws.Cells.Item(6, 14).Value = 4.1
$ws.Cells.Item(6, 15).Value = 1.25
$ws.Cells.Item(6, 16).Value = 2.08
$ws.Cells.Item(6, 17).Value = 1.74
$ws.Cells.Item(6, 18).Value = 1.42
$ws.Cells.Item(6, 19).Value = 2.86
$ws.Cells.Item(6, 20).Value = 1.64
$ws.Cells.Item(6, 21).Value = 2.24
$ws.Cells.Item(6, 22).Value = 1.4
$ws.Cells.Item(6, 23).Value = 1.72
$ws.Cells.Item(6, 24).Value = 18.5
$ws.Cells.Item(6, 25).Value = 16
$ws.Cells.Item(6, 26).Value = 30
$ws.Cells.Item(6, 27).Value = 60
$ws.Cells.Item(6, 28).Value = 12.5
$ws.Cells.Item(6, 29).Value = 10
$ws.Cells.Item(6, 30).Value = 17
$ws.Cells.Item(6, 31).Value = 42
$ws.Cells.Item(6, 32).Value = 16.5
$ws.Cells.Item(6, 33).Value = 13
$ws.Cells.Item(6, 34).Value = 19
$ws.Cells.Item(6, 35).Value = 50
$ws.Cells.Item(6, 36).Value = 34
$ws.Cells.Item(6, 37).Value = 26
$ws.Cells.Item(6, 38).Value = 38
$ws.Cells.Item(6, 39).Value = 90
$ws.Cells.Item(6, 40).Value = 15.5
$ws.Cells.Item(6, 41).Value = 32

# Row 7: Scottish League Two | East Kilbride vs Spartans
$ws.Cells.Item(7, 1).Value = "Scottish League Two"
$ws.Cells.Item(7, 2).Value = "2026-01-27"
$ws.Cells.Item(7, 3).Value = "16:45:00"
$ws.Cells.Item(7, 4).Value = "East Kilbride"
$ws.Cells.Item(7, 5).Value = "Spartans"
$ws.Cells.Item(7, 6).Value = 1.04
$ws.Cells.Item(7, 7).Value = 1000
$ws.Cells.Item(7, 8).Value = 1.04
$ws.Cells.Item(7, 9).Value = 1000
$ws.Cells.Item(7, 10).Value = 1.02
$ws.Cells.Item(7, 11).Value = 950
$ws.Cells.Item(7, 12).Value = 1.01
$ws.Cells.Item(7, 13).Value = 1.01
$ws.Cells.Item(7, 14).Value = 1.08
$ws.Cells.Item(7, 15).Value = 1.24
$ws.Cells.Item(7, 16).Value = 1.08
$ws.Cells.Item(7, 17).Value = 1.24
$ws.Cells.Item(7, 18).Value = 1.08
$ws.Cells.Item(7, 19).Value = 1.24
$ws.Cells.Item(7, 20).Value = 1.01
$ws.Cells.Item(7, 21).Value = 1.01
$ws.Cells.Item(7, 22).Value = 1.01
$ws.Cells.Item(7, 23).Value = 1.01
$ws.Cells.Item(7, 24).Value = 1000
$ws.Cells.Item(7, 25).Value = 1000
$ws.Cells.Item(7, 26).Value = 1000
$ws.Cells.Item(7, 27).Value = 1000
$ws.Cells.Item(7, 28).Value = 1000
$ws.Cells.Item(7, 29).Value = 1000
$ws.Cells.Item(7, 30).Value = 1000
$ws.Cells.Item(7, 31).Value = 1000
$ws.Cells.Item(7, 32).Value = 1000
$ws.Cells.Item(7, 33).Value = 1000
$ws.Cells.Item(7, 34).Value = 1000
$ws.Cells.Item(7, 35).Value = 1000
$ws.Cells.Item(7, 36).Value = 1000
$ws.Cells.Item(7, 37).Value = 1000
$ws.Cells.Item(7, 38).Value = 1000
$ws.Cells.Item(7, 39).Value = 1000
$ws.Cells.Item(7, 40).Value = 1000
$ws.Cells.Item(7, 41).Value = 1000

# Row 8: Scottish League Two | Stranraer vs Clyde
$ws.Cells.Item(8, 1).Value = "Scottish League Two"
$ws.Cells.Item(8, 2).Value = "2026-01-27"
$ws.Cells.Item(8, 3).Value = "16:45:00"
$ws.Cells.Item(8, 4).Value = "Stranraer"
$ws.Cells.Item(8, 5).Value = "Clyde"
$ws.Cells.Item(8, 6).Value = 1.45
$ws.Cells.Item(8, 7).Value = 3
$ws.Cells.Item(8, 8).Value = 1.53
$ws.Cells.Item(8, 9).Value = 3.4
$ws.Cells.Item(8, 10).Value = 2.96
$ws.Cells.Item(8, 11).Value = 950
$ws.Cells.Item(8, 12).Value = 1.3
$ws.Cells.Item(8, 13).Value = 1.01
$ws.Cells.Item(8, 14).Value = 1.08
$ws.Cells.Item(8, 15).Value = 1.01
$ws.Cells.Item(8, 16).Value = 1.08
$ws.Cells.Item(8, 17).Value = 1.01
$ws.Cells.Item(8, 18).Value = 1.08
$ws.Cells.Item(8, 19).Value = 1.01
$ws.Cells.Item(8, 20).Value = 1.01
$ws.Cells.Item(8, 21).Value = 1.01
$ws.Cells.Item(8, 22).Value = 1.41
$ws.Cells.Item(8, 23).Value = 1.5
$ws.Cells.Item(8, 24).Value = 1000
$ws.Cells.Item(8, 25).Value = 1000
$ws.Cells.Item(8, 26).Value = 1000
$ws.Cells.Item(8, 27).Value = 1000
$ws.Cells.Item(8, 28).Value = 1000
$ws.Cells.Item(8, 29).Value = 1000
$ws.Cells.Item(8, 30).Value = 1000
$ws.Cells.Item(8, 31).Value = 1000
$ws.Cells.Item(8, 32).Value = 1000
$ws.Cells.Item(8, 33).Value = 1000
$ws.Cells.Item(8, 34).Value = 1000
$ws.Cells.Item(8, 35).Value = 1000
$ws.Cells.Item(8, 36).Value = 1000
$ws.Cells.Item(8, 37).Value = 1000
$ws.Cells.Item(8, 38).Value = 1000
$ws.Cells.Item(8, 39).Value = 1000
$ws.Cells.Item(8, 40).Value = 1000
$ws.Cells.Item(8, 41).Value = 1000

# Row 9: Argentinian Primera Division | Velez Sarsfield vs Talleres
$ws.Cells.Item(9, 1).Value = "Argentinian Primera Division"
$ws.Cells.Item(9, 2).Value = "2026-01-27"
$ws.Cells.Item(9, 3).Value = "17:45:00"
$ws.Cells.Item(9, 4).Value = "Velez Sarsfield"
$ws.Cells.Item(9, 5).Value = "Talleres"
$ws.Cells.Item(9, 6).Value = 2.3
$ws.Cells.Item(9, 7).Value = 2.48
$ws.Cells.Item(9, 8).Value = 3.85
$ws.Cells.Item(9, 9).Value = 4.3
$ws.Cells.Item(9, 10).Value = 2.88
$ws.Cells.Item(9, 11).Value = 3.15
$ws.Cells.Item(9, 12).Value = 0
$ws.Cells.Item(9, 13).Value = 0
$ws.Cells.Item(9, 14).Value = 0
$ws.Cells.Item(9, 15).Value = 0
$ws.Cells.Item(9, 16).Value = 1.07
$ws.Cells.Item(9, 17).Value = 1.01
$ws.Cells.Item(9, 18).Value = 0
$ws.Cells.Item(9, 19).Value = 0
$ws.Cells.Item(9, 20).Value = 0
$ws.Cells.Item(9, 21).Value = 0
$ws.Cells.Item(9, 22).Value = 0
$ws.Cells.Item(9, 23).Value = 0
$ws.Cells.Item(9, 24).Value = 0
$ws.Cells.Item(9, 25).Value = 0
$ws.Cells.Item(9, 26).Value = 0
$ws.Cells.Item(9, 27).Value = 0
$ws.Cells.Item(9, 28).Value = 0
$ws.Cells.Item(9, 29).Value = 0
$ws.Cells.Item(9, 30).Value = 0
$ws.Cells.Item(9, 31).Value = 0
$ws.Cells.Item(9, 32).Value = 0
$ws.Cells.Item(9, 33).Value = 0
$ws.Cells.Item(9, 34).Value = 0
$ws.Cells.Item(9, 35).Value = 0
$ws.Cells.Item(9, 36).Value = 0
$ws.Cells.Item(9, 37).Value = 0
$ws.Cells.Item(9, 38).Value = 0
$ws.Cells.Item(9, 39).Value = 0
$ws.Cells.Item(9, 40).Value = 0
$ws.Cells.Item(9, 41).Value = 0

# Row 10: Colombian Primera A | Cucuta Deportivo vs Atletico Bucaramanga
$ws.Cells.Item(10, 1).Value = "Colombian Primera A"
$ws.Cells.Item(10, 2).Value = "2026-01-27"
$ws.Cells.Item(10, 3).Value = "18:00:00"
$ws.Cells.Item(10, 4).Value = "Cucuta Deportivo"
$ws.Cells.Item(10, 5).Value = "Atletico Bucaramanga"
$ws.Cells.Item(10, 6).Value = 2.56
$ws.Cells.Item(10, 7).Value = 3.65
$ws.Cells.Item(10, 8).Value = 2.64
$ws.Cells.Item(10, 9).Value = 3.25
$ws.Cells.Item(10, 10).Value = 2.68
$ws.Cells.Item(10, 11).Value = 3.65
$ws.Cells.Item(10, 12).Value = 0
$ws.Cells.Item(10, 13).Value = 0
$ws.Cells.Item(10, 14).Value = 0
$ws.Cells.Item(10, 15).Value = 0
$ws.Cells.Item(10, 16).Value = 1.4
$ws.Cells.Item(10, 17).Value = 1.01
$ws.Cells.Item(10, 18).Value = 0
$ws.Cells.Item(10, 19).Value = 0
$ws.Cells.Item(10, 20).Value = 0
$ws.Cells.Item(10, 21).Value = 0
$ws.Cells.Item(10, 22).Value = 0
$ws.Cells.Item(10, 23).Value = 0
$ws.Cells.Item(10, 24).Value = 0
$ws.Cells.Item(10, 25).Value = 0
$ws.Cells.Item(10, 26).Value = 0
$ws.Cells.Item(10, 27).Value = 0
$ws.Cells.Item(10, 28).Value = 0
$ws.Cells.Item(10, 29).Value = 0
$ws.Cells.Item(10, 30).Value = 0
$ws.Cells.Item(10, 31).Value = 0
$ws.Cells.Item(10, 32).Value = 0
$ws.Cells.Item(10, 33).Value = 0
$ws.Cells.Item(10, 34).Value = 0
$ws.Cells.Item(10, 35).Value = 0
$ws.Cells.Item(10, 36).Value = 0
$ws.Cells.Item(10, 37).Value = 0
$ws.Cells.Item(10, 38).Value = 0
$ws.Cells.Item(10, 39).Value = 0
$ws.Cells.Item(10, 40).Value = 0
$ws.Cells.Item(10, 41).Value = 0

# Row 11: Paraguayan Primera Division | Sportivo Luqueno vs Nacional (Par)
$ws.Cells.Item(11, 1).Value = "Paraguayan Primera Division"
$ws.Cells.Item(11, 2).Value = "2026-01-27"
$ws.Cells.Item(11, 3).Value = "18:00:00"
$ws.Cells.Item(11, 4).Value = "Sportivo Luqueno"
$ws.Cells.Item(11, 5).Value = "Nacional (Par)"
$ws.Cells.Item(11, 6).Value = 0
$ws.Cells.Item(11, 7).Value = 0
$ws.Cells.Item(11, 8).Value = 0
$ws.Cells.Item(11, 9).Value = 0
$ws.Cells.Item(11, 10).Value = 0
$ws.Cells.Item(11, 11).Value = 0
$ws.Cells.Item(11, 12).Value = 0
$ws.Cells.Item(11, 13).Value = 0
$ws.Cells.Item(11, 14).Value = 0
$ws.Cells.Item(11, 15).Value = 0
$ws.Cells.Item(11, 16).Value = 1.07
$ws.Cells.Item(11, 17).Value = 1.01
$ws.Cells.Item(11, 18).Value = 0
$ws.Cells.Item(11, 19).Value = 0
$ws.Cells.Item(11, 20).Value = 0
$ws.Cells.Item(11, 21).Value = 0
$ws.Cells.Item(11, 22).Value = 0
$ws.Cells.Item(11, 23).Value = 0
$ws.Cells.Item(11, 24).Value = 0
$ws.Cells.Item(11, 25).Value = 0
$ws.Cells.Item(11, 26).Value = 0
$ws.Cells.Item(11, 27).Value = 0
$ws.Cells.Item(11, 28).Value = 0
$ws.Cells.Item(11, 29).Value = 0
$ws.Cells.Item(11, 30).Value = 0
$ws.Cells.Item(11, 31).Value = 0
$ws.Cells.Item(11, 32).Value = 0
$ws.Cells.Item(11, 33).Value = 0
$ws.Cells.Item(11, 34).Value = 0
$ws.Cells.Item(11, 35).Value = 0
$ws.Cells.Item(11, 36).Value = 0
$ws.Cells.Item(11, 37).Value = 0
$ws.Cells.Item(11, 38).Value = 0
$ws.Cells.Item(11, 39).Value = 0
$ws.Cells.Item(11, 40).Value = 0
$ws.Cells.Item(11, 41).Value = 0

# Row 12: Argentinian Primera Division | CA Huracan vs Independiente Rivadavia
$ws.Cells.Item(12, 1).Value = "Argentinian Primera Division"
$ws.Cells.Item(12, 2).Value = "2026-01-27"
$ws.Cells.Item(12, 3).Value = "20:00:00"
$ws.Cells.Item(12, 4).Value = "CA Huracan"
$ws.Cells.Item(12, 5).Value = "Independiente Rivadavia"
$ws.Cells.Item(12, 6).Value = 2.02
$ws.Cells.Item(12, 7).Value = 2.18
$ws.Cells.Item(12, 8).Value = 4.7
$ws.Cells.Item(12, 9).Value = 5.4
$ws.Cells.Item(12, 10).Value = 2.88
$ws.Cells.Item(12, 11).Value = 3.3
$ws.Cells.Item(12, 12).Value = 0
$ws.Cells.Item(12, 13).Value = 0
$ws.Cells.Item(12, 14).Value = 0
$ws.Cells.Item(12, 15).Value = 0
$ws.Cells.Item(12, 16).Value = 1.42
$ws.Cells.Item(12, 17).Value = 2.96
$ws.Cells.Item(12, 18).Value = 0
$ws.Cells.Item(12, 19).Value = 0
$ws.Cells.Item(12, 20).Value = 0
$ws.Cells.Item(12, 21).Value = 0
$ws.Cells.Item(12, 22).Value = 0
$ws.Cells.Item(12, 23).Value = 0
$ws.Cells.Item(12, 24).Value = 0
$ws.Cells.Item(12, 25).Value = 0
$ws.Cells.Item(12, 26).Value = 0
$ws.Cells.Item(12, 27).Value = 0
$ws.Cells.Item(12, 28).Value = 0
$ws.Cells.Item(12, 29).Value = 0
$ws.Cells.Item(12, 30).Value = 0
$ws.Cells.Item(12, 31).Value = 0
$ws.Cells.Item(12, 32).Value = 0
$ws.Cells.Item(12, 33).Value = 0
$ws.Cells.Item(12, 34).Value = 0
$ws.Cells.Item(12, 35).Value = 0
$ws.Cells.Item(12, 36).Value = 0
$ws.Cells.Item(12, 37).Value = 0
$ws.Cells.Item(12, 38).Value = 0
$ws.Cells.Item(12, 39).Value = 0
$ws.Cells.Item(12, 40).Value = 0
$ws.Cells.Item(12, 41).Value = 0

# Row 13: Argentinian Primera Division | Gimnasia Mendoza vs San Lorenzo
$ws.Cells.Item(13, 1).Value = "Argentinian Primera Division"
$ws.Cells.Item(13, 2).Value = "2026-01-27"
$ws.Cells.Item(13, 3).Value = "20:00:00"
$ws.Cells.Item(13, 4).Value = "Gimnasia Mendoza"
$ws.Cells.Item(13, 5).Value = "San Lorenzo"
$ws.Cells.Item(13, 6).Value = 2.88
$ws.Cells.Item(13, 7).Value = 3.15
$ws.Cells.Item(13, 8).Value = 3.05
$ws.Cells.Item(13, 9).Value = 3.35
$ws.Cells.Item(13, 10).Value = 2.7
$ws.Cells.Item(13, 11).Value = 2.88
$ws.Cells.Item(13, 12).Value = 0
$ws.Cells.Item(13, 13).Value = 0
$ws.Cells.Item(13, 14).Value = 0
$ws.Cells.Item(13, 15).Value = 0
$ws.Cells.Item(13, 16).Value = 1.35
$ws.Cells.Item(13, 17).Value = 3.5
$ws.Cells.Item(13, 18).Value = 0
$ws.Cells.Item(13, 19).Value = 0
$ws.Cells.Item(13, 20).Value = 0
$ws.Cells.Item(13, 21).Value = 0
$ws.Cells.Item(13, 22).Value = 0
$ws.Cells.Item(13, 23).Value = 0
$ws.Cells.Item(13, 24).Value = 0
$ws.Cells.Item(13, 25).Value = 0
$ws.Cells.Item(13, 26).Value = 0
$ws.Cells.Item(13, 27).Value = 0
$ws.Cells.Item(13, 28).Value = 0
$ws.Cells.Item(13, 29).Value = 0
$ws.Cells.Item(13, 30).Value = 0
$ws.Cells.Item(13, 31).Value = 0
$ws.Cells.Item(13, 32).Value = 0
$ws.Cells.Item(13, 33).Value = 0
$ws.Cells.Item(13, 34).Value = 0
$ws.Cells.Item(13, 35).Value = 0
$ws.Cells.Item(13, 36).Value = 0
$ws.Cells.Item(13, 37).Value = 0
$ws.Cells.Item(13, 38).Value = 0
$ws.Cells.Item(13, 39).Value = 0
$ws.Cells.Item(13, 40).Value = 0
$ws.Cells.Item(13, 41).Value = 0

# Row 14: Colombian Primera A | Fortaleza FC vs Llaneros FC
$ws.Cells.Item(14, 1).Value = "Colombian Primera A"
$ws.Cells.Item(14, 2).Value = "2026-01-27"
$ws.Cells.Item(14, 3).Value = "20:20:00"
$ws.Cells.Item(14, 4).Value = "Fortaleza FC"
$ws.Cells.Item(14, 5).Value = "Llaneros FC"
$ws.Cells.Item(14, 6).Value = 1.93
$ws.Cells.Item(14, 7).Value = 2.14
$ws.Cells.Item(14, 8).Value = 4.4
$ws.Cells.Item(14, 9).Value = 5.9
$ws.Cells.Item(14, 10).Value = 2.82
$ws.Cells.Item(14, 11).Value = 3.65
$ws.Cells.Item(14, 12).Value = 0
$ws.Cells.Item(14, 13).Value = 0
$ws.Cells.Item(14, 14).Value = 0
$ws.Cells.Item(14, 15).Value = 0
$ws.Cells.Item(14, 16).Value = 1.57
$ws.Cells.Item(14, 17).Value = 2.1
$ws.Cells.Item(14, 18).Value = 0
$ws.Cells.Item(14, 19).Value = 0
$ws.Cells.Item(14, 20).Value = 0
$ws.Cells.Item(14, 21).Value = 0
$ws.Cells.Item(14, 22).Value = 0
$ws.Cells.Item(14, 23).Value = 0
$ws.Cells.Item(14, 24).Value = 0
$ws.Cells.Item(14, 25).Value = 0
$ws.Cells.Item(14, 26).Value = 0
$ws.Cells.Item(14, 27).Value = 0
$ws.Cells.Item(14, 28).Value = 0
$ws.Cells.Item(14, 29).Value = 0
$ws.Cells.Item(14, 30).Value = 0
$ws.Cells.Item(14, 31).Value = 0
$ws.Cells.Item(14, 32).Value = 0
$ws.Cells.Item(14, 33).Value = 0
$ws.Cells.Item(14, 34).Value = 0
$ws.Cells.Item(14, 35).Value = 0
$ws.Cells.Item(14, 36).Value = 0
$ws.Cells.Item(14, 37).Value = 0
$ws.Cells.Item(14, 38).Value = 0
$ws.Cells.Item(14, 39).Value = 0
$ws.Cells.Item(14, 40).Value = 0
$ws.Cells.Item(14, 41).Value = 0
